$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing respondent row 33: Python5 flag 0 -> 1 ---
$ws.Range("L33").Value = 1

# --- Append new respondent as row 37 (#36) ---
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "Ana Salet Hidalgo Flores"
$ws.Range("C37").Value = "anitasalet2203@gmail.com"
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 8
$ws.Range("H37").Value = 1
$ws.Range("I37").Value = 3
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0

# The new row was pasted in from another workbook, which (as in real Excel)
# brings along a new "Normal 2" cell style for the name/email columns.
$normal2 = $wb.Styles.Add("Normal 2")
$normal2.Font.Name = "Arial"
$normal2.Font.Size = 10
$normal2.Font.Color = 0
$ws.Range("B37").Style = "Normal 2"
$ws.Range("C37").Style = "Normal 2"

# Re-select near the new row, like after scrolling down to review/enter it.
$ws.Range("M33").Select()

Write-Host "done"
